$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01037833333333333
$ws.Range("H2").Value = 0.031135
$ws.Range("I2").Value = 0.02114284782989566
$ws.Range("J2").Value = 0.02114284782989566
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 0.01497632937666667
$ws.Range("R2").Value = 0.13478696439
$ws.Range("S2").Value = 0.0006124382783924709
$ws.Range("T2").Value = 0.0006124382783924709
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01037833333333333
$ws.Range("H3").Value = 0.031135
$ws.Range("I3").Value = 0.02114284782989566
$ws.Range("J3").Value = 0.02114284782989566
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 0.3030838062561111
$ws.Range("R3").Value = 2.727754256305
$ws.Range("S3").Value = 0.01239423358311875
$ws.Range("T3").Value = 0.01239423358311875
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01037833333333333
$ws.Range("H4").Value = 0.031135
$ws.Range("I4").Value = 0.02114284782989566
$ws.Range("J4").Value = 0.02114284782989566
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 0.198958908135
$ws.Range("R4").Value = 1.790630173215
$ws.Range("S4").Value = 0.008136175968384436
$ws.Range("T4").Value = 0.008136175968384436
$ws.Range("G5").Value = 0.4265683333333333
$ws.Range("I5").Value = 0.8690094132698448
$ws.Range("J5").Value = 0.8690094132698448
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 0.6155543145966665
$ws.Range("R5").Value = 5.539988831369999
$ws.Range("S5").Value = 0.02517232462020995
$ws.Range("T5").Value = 0.02517232462020995
$ws.Range("G6").Value = 0.4265683333333333
$ws.Range("I6").Value = 0.8690094132698448
$ws.Range("J6").Value = 0.8690094132698448
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("S6").Value = 0.5094254918093778
$ws.Range("T6").Value = 0.5094254918093778
$ws.Range("G7").Value = 0.4265683333333333
$ws.Range("I7").Value = 0.8690094132698448
$ws.Range("J7").Value = 0.8690094132698448
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 8.177572170704998
$ws.Range("R7").Value = 73.59814953634499
$ws.Range("S7").Value = 0.3344115968402571
$ws.Range("T7").Value = 0.334411596840257
$ws.Range("G8").Value = 0.05392066666666667
$ws.Range("H8").Value = 0.161762
$ws.Range("I8").Value = 0.1098477389002595
$ws.Range("J8").Value = 0.1098477389002595
$ws.Range("M8").Value = 1.443038
$ws.Range("N8").Value = 4.329114
$ws.Range("O8").Value = 0.0289666880885598
$ws.Range("P8").Value = 0.0289666880885598
$ws.Range("Q8").Value = 0.07780957098533332
$ws.Range("R8").Value = 0.700286138868
$ws.Range("S8").Value = 0.003181925189957375
$ws.Range("T8").Value = 0.003181925189957375
$ws.Range("G9").Value = 0.05392066666666667
$ws.Range("H9").Value = 0.161762
$ws.Range("I9").Value = 0.1098477389002595
$ws.Range("J9").Value = 0.1098477389002595
$ws.Range("N9").Value = 87.61054300000001
$ws.Range("O9").Value = 0.5862140087672342
$ws.Range("P9").Value = 0.5862140087672342
$ws.Range("Q9").Value = 1.574672961862889
$ws.Range("R9").Value = 14.172056656766
$ws.Range("S9").Value = 0.0643942833747376
$ws.Range("T9").Value = 0.0643942833747376
$ws.Range("G10").Value = 0.05392066666666667
$ws.Range("H10").Value = 0.161762
$ws.Range("I10").Value = 0.1098477389002595
$ws.Range("J10").Value = 0.1098477389002595
$ws.Range("M10").Value = 19.170603
$ws.Range("N10").Value = 57.511809
$ws.Range("O10").Value = 0.384819303144206
$ws.Range("P10").Value = 0.384819303144206
$ws.Range("Q10").Value = 1.033691694162
$ws.Range("R10").Value = 9.303225247458
$ws.Range("S10").Value = 0.04227153033556457
$ws.Range("T10").Value = 0.04227153033556456
